$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Extend formatting from the last populated row (90) down through the
# three new rows (91-93) so the new rows pick up the same number formats,
# fonts and borders as the rest of the data table.
$ws.Range("A90:I90").Copy()
$ws.Range("A91:I93").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 91
$ws.Range("A91").Value = 88
$ws.Range("B91").Value = 26878451
$ws.Range("C91").Value = 46426
$ws.Range("D91").Value = 44000.67
$ws.Range("E91").Value = 44321
$ws.Range("F91").Formula = '=IF(B91="","",C91-D91)'
$ws.Range("G91").Formula = '=IF(B91="","",F91/D91*100)'
$ws.Range("H91").Formula = '=IF(B91="","",D91*1.04)'
$ws.Range("I91").Formula = '=IF(B91="","",C91-H91)'

# Row 92
$ws.Range("A92").Value = 89
$ws.Range("B92").Value = 26883517
$ws.Range("C92").Value = 211026
$ws.Range("D92").Value = 200000.32
$ws.Range("E92").Value = 44322
$ws.Range("F92").Formula = '=IF(B92="","",C92-D92)'
$ws.Range("G92").Formula = '=IF(B92="","",F92/D92*100)'
$ws.Range("H92").Formula = '=IF(B92="","",D92*1.04)'
$ws.Range("I92").Formula = '=IF(B92="","",C92-H92)'

# Row 93
$ws.Range("A93").Value = 90
$ws.Range("B93").Value = 26903291
$ws.Range("C93").Value = 113421
$ws.Range("D93").Value = 107494.91
$ws.Range("E93").Value = 44324
$ws.Range("F93").Formula = '=IF(B93="","",C93-D93)'
$ws.Range("G93").Formula = '=IF(B93="","",F93/D93*100)'
$ws.Range("H93").Formula = '=IF(B93="","",D93*1.04)'
$ws.Range("I93").Formula = '=IF(B93="","",C93-H93)'

$ws.Range("E94").Select()
